# The sheet holds quarterly data in blocks of 4 rows per year
# (row labels "<year>年A/B/C/D"). For every year block, the "B" row and
# "C" row need to swap places (their A:E contents trade rows), while the
# "A" and "D" rows stay put. After that, the no-longer-needed "F" and "G"
# columns (machine-made paper & paperboard production-sales ratio /
# sales volume, duplicated from columns B/E) are deleted entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "C", "D", "E")

for ($yearStart = 2; $yearStart -le 61; $yearStart += 4) {
    $rowB = $yearStart + 1
    $rowC = $yearStart + 2

    foreach ($col in $cols) {
        $cellB = $ws.Range($col + $rowB)
        $cellC = $ws.Range($col + $rowC)

        $valB = $cellB.Value2
        $valC = $cellC.Value2

        $cellB.Value2 = $valC
        $cellC.Value2 = $valB
    }
}

# Remove the now-redundant F (产销率) and G (销售量) columns.
$ws.Range("F1:G1").EntireColumn.Delete()
